{"js": "// Remove the \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" paragraphs, along with the\n// blank spacer paragraph that precedes them (directly after the\n// \"LOB1019: F\u00edsica II (Requisito fraco)\" paragraph), matching the diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (\n    t === \"Ver no Jupiter Salvar em pdf Salvar em docx\" ||\n    t === \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n  ) {\n    targets.push(i);\n  }\n}\n\n// The blank paragraph immediately before the first target paragraph (the\n// spacer inserted after the \"LOB1019\" line) is removed too.\nif (targets.length > 0) {\n  const firstIdx = targets[0];\n  const prev = paragraphs.items[firstIdx - 1];\n  if (prev && prev.text === \"\") {\n    targets.unshift(firstIdx - 1);\n  }\n}\n\nfor (const idx of targets) {\n  paragraphs.items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" paragraphs, along with the\n# blank spacer paragraph that precedes them (directly after the\n# \"LOB1019: F\u00edsica II (Requisito fraco)\" paragraph), matching the diff.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $jupiterText) {\n        $jupiterIdx = $i\n    } elseif ($t -eq $copyrightText) {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -gt 0 -and $copyrightIdx -gt 0) {\n    $blankIdx = $jupiterIdx - 1\n    $blankText = $d.Paragraphs.Item($blankIdx).Range.Text.TrimEnd([char]13, [char]7)\n\n    # Delete from the last target paragraph back to the first so earlier\n    # indices stay valid as later ones are removed.\n    $d.Paragraphs.Item($copyrightIdx).Range.Delete()\n    $d.Paragraphs.Item($jupiterIdx).Range.Delete()\n    if ($blankText -eq \"\") {\n        $d.Paragraphs.Item($blankIdx).Range.Delete()\n    }\n}\n"}
